$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column C for rows 2-11 from 45208 to 45212
$ws.Range("C2:C11").Value = 45212
